$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values can look like plain numbers (e.g. "1.000", "6.050", "0.2180"),
# which Excel would otherwise auto-convert to numeric and mangle (trailing
# zeros / float rounding lost). Pre-format the affected range as Text so the
# values are stored verbatim as strings, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '28.175.54'
$ws.Cells.Item(2, 5).Value = '  -0.27%  '
$ws.Cells.Item(3, 4).Value = '1.827.67'
$ws.Cells.Item(3, 5).Value = '  +1.24%  '
$ws.Cells.Item(4, 4).Value = '0.9992'
$ws.Cells.Item(4, 5).Value = '  -0.45%  '
$ws.Cells.Item(5, 4).Value = '310.45'
$ws.Cells.Item(5, 5).Value = '  -1.33%  '
$ws.Cells.Item(6, 4).Value = '0.9995'
$ws.Cells.Item(6, 5).Value = '  -0.28%  '
$ws.Cells.Item(7, 4).Value = '0.5135'
$ws.Cells.Item(7, 5).Value = '  -2.60%  '
$ws.Cells.Item(8, 4).Value = '0.3925'
$ws.Cells.Item(8, 5).Value = '  +2.55%  '
$ws.Cells.Item(9, 4).Value = '0.09609'
$ws.Cells.Item(9, 5).Value = '  +19.78%  '
$ws.Cells.Item(10, 4).Value = '1.109'
$ws.Cells.Item(10, 5).Value = '  +0.65%  '
$ws.Cells.Item(11, 4).Value = '40.93'
$ws.Cells.Item(11, 5).Value = '  -1.20%  '
$ws.Cells.Item(12, 4).Value = '6.466'
$ws.Cells.Item(12, 5).Value = '  +2.12%  '
$ws.Cells.Item(13, 4).Value = '0.9991'
$ws.Cells.Item(13, 5).Value = '  -0.44%  '
$ws.Cells.Item(14, 4).Value = '20.59'
$ws.Cells.Item(14, 5).Value = '  -0.16%  '
$ws.Cells.Item(15, 4).Value = '1.816.81'
$ws.Cells.Item(15, 5).Value = '  +0.49%  '
$ws.Cells.Item(16, 4).Value = '7.367'
$ws.Cells.Item(16, 5).Value = '  +0.44%  '
$ws.Cells.Item(17, 4).Value = '0.00001136'
$ws.Cells.Item(17, 5).Value = '  +3.58%  '
$ws.Cells.Item(18, 4).Value = '92.67'
$ws.Cells.Item(18, 5).Value = '  +0.49%  '
$ws.Cells.Item(19, 4).Value = '0.06592'
$ws.Cells.Item(19, 5).Value = '  -0.25%  '
$ws.Cells.Item(20, 4).Value = '1.000'
$ws.Cells.Item(20, 5).Value = '  -0.19%  '
$ws.Cells.Item(21, 4).Value = '17.32'
$ws.Cells.Item(21, 5).Value = '  -0.45%  '
$ws.Cells.Item(22, 4).Value = '6.050'
$ws.Cells.Item(22, 5).Value = '  +1.30%  '
$ws.Cells.Item(23, 4).Value = '28.241.56'
$ws.Cells.Item(23, 5).Value = '  -0.26%  '
$ws.Cells.Item(24, 4).Value = '11.17'
$ws.Cells.Item(24, 5).Value = '  +0.07%  '
$ws.Cells.Item(25, 4).Value = '2.228'
$ws.Cells.Item(25, 5).Value = '  -1.30%  '
$ws.Cells.Item(26, 4).Value = '157.42'
$ws.Cells.Item(26, 5).Value = '  -2.19%  '
$ws.Cells.Item(27, 4).Value = '2.445'
$ws.Cells.Item(27, 5).Value = '  +3.42%  '
$ws.Cells.Item(28, 4).Value = '20.61'
$ws.Cells.Item(28, 5).Value = '  +0.63%  '
$ws.Cells.Item(29, 4).Value = '2.016.80'
$ws.Cells.Item(29, 5).Value = '  +0.31%  '
$ws.Cells.Item(30, 4).Value = '128.64'
$ws.Cells.Item(30, 5).Value = '  +4.25%  '
$ws.Cells.Item(31, 4).Value = '0.1092'
$ws.Cells.Item(31, 5).Value = '  +0.69%  '
$ws.Cells.Item(32, 4).Value = '1.061'
$ws.Cells.Item(32, 5).Value = '  +0.50%  '
$ws.Cells.Item(33, 4).Value = '5.659'
$ws.Cells.Item(33, 5).Value = '  +1.70%  '
$ws.Cells.Item(34, 4).Value = '3.631'
$ws.Cells.Item(34, 5).Value = '  -1.46%  '
$ws.Cells.Item(35, 4).Value = '0.06933'
$ws.Cells.Item(35, 5).Value = '  -4.84%  '
$ws.Cells.Item(36, 4).Value = '9.108'
$ws.Cells.Item(36, 5).Value = '  +5.06%  '
$ws.Cells.Item(37, 4).Value = '0.02341'
$ws.Cells.Item(37, 5).Value = '  +0.95%  '
$ws.Cells.Item(38, 4).Value = '0.2180'
$ws.Cells.Item(38, 5).Value = '  +0.92%  '
$ws.Cells.Item(39, 4).Value = '11.59'
$ws.Cells.Item(39, 5).Value = '  -7.39%  '
$ws.Cells.Item(40, 4).Value = '5.035'
$ws.Cells.Item(40, 5).Value = '  -1.70%  '
$ws.Cells.Item(41, 4).Value = '0.6209'
$ws.Cells.Item(41, 5).Value = '  -0.03%  '
$ws.Cells.Item(42, 4).Value = '0.9988'
$ws.Cells.Item(42, 5).Value = '  -0.26%  '
$ws.Cells.Item(43, 4).Value = '1.156'
$ws.Cells.Item(43, 5).Value = '  -1.17%  '
$ws.Cells.Item(44, 4).Value = '13.31'
$ws.Cells.Item(44, 5).Value = '  +0.75%  '
$ws.Cells.Item(45, 4).Value = '0.5979'
$ws.Cells.Item(45, 5).Value = '  -0.93%  '
$ws.Cells.Item(46, 4).Value = '1.291'
$ws.Cells.Item(46, 5).Value = '  -5.69%  '
$ws.Cells.Item(47, 4).Value = '3.709'
$ws.Cells.Item(47, 5).Value = '  -1.66%  '
$ws.Cells.Item(48, 4).Value = '125.48'
$ws.Cells.Item(48, 5).Value = '  -1.35%  '
$ws.Cells.Item(49, 4).Value = '1.964'
$ws.Cells.Item(49, 5).Value = '  +1.63%  '
$ws.Cells.Item(50, 4).Value = '1.188'
$ws.Cells.Item(50, 5).Value = '  -2.67%  '
$ws.Cells.Item(51, 4).Value = '0.06784'
$ws.Cells.Item(51, 5).Value = '  -0.49%  '
